$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("Z3").Value = "Convalescent"
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Value = "12"
$ws.Range("AU3").Value = "month"
$ws.Range("AV3").Value = "0 - 9"
$ws.Range("BA3").Value = "Ageusia (complete loss of taste); Cough; Headache; Hypotension (low blood pressure)"
$ws.Range("BF3").Value = "New York"

# Row 4
$ws.Range("Z4").Value = "Convalescent"
$ws.Range("AU4").Value = "year"
$ws.Range("BA4").Value = "Confusion; Cough; Pharyngitis (sore throat)"

# Row 5
$ws.Range("S5").Value = "Canada"
$ws.Range("T5").Value = "Alberta"
$ws.Range("AU5").Value = "year"
$ws.Range("BA5").Value = "Coma; Dyspnea (breathing difficulty); Rhinorrhea (runny nose)"
$ws.Range("BF5").Value = "Dublin"
$ws.Range("BK5").Value = "Missing"

# Row 6
$ws.Range("AN6").Value = "Human"
$ws.Range("AO6").Value = "Homo sapiens"
$ws.Range("AU6").Value = "year"
$ws.Range("BA6").Value = "Confusion; Irritability; Cough"
$ws.Range("BF6").Value = "Provence-Alpes-Cote d'Azur"

# Row 7
$ws.Range("AU7").Value = "year"
$ws.Range("BA7").Value = "Cognitive impairment; Fever"
$ws.Range("BF7").Value = "Gauteng"
